$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows (Y/N rows no longer needed) - rows 3,4,5
$ws.Cells.Item(4, 1).ClearContents()
$ws.Cells.Item(5, 1).ClearContents()

# Row 2: yousuf / 793206 / TRUE / FAILED
$ws.Cells.Item(2, 1).Value = "yousuf"
$ws.Cells.Item(2, 2).Value = 793206
$ws.Cells.Item(2, 3).Value = $true

# Row 3: riyak / 110223 / FALSE / FAILED
$ws.Cells.Item(3, 1).Value = "riyak"
$ws.Cells.Item(3, 2).Value = 110223
$ws.Cells.Item(3, 3).Value = $false

# E2 typed as PASSED first, then corrected to FAILED (leaves PASSED as orphan shared string)
$ws.Cells.Item(2, 5).Value = "PASSED"
$ws.Cells.Item(2, 5).Value = "FAILED"
$ws.Cells.Item(3, 5).Value = "FAILED"

$ws.Range("C3").Select()
